# TaskTwo_Automation_one.xlsx
#
# Add a new "Report utility" worksheet (FindNewUser) that holds the
# "find/search a user" record-table rows which used to live at the bottom
# of AddUser (rows 23:25), give it its own header row, and remove the
# now-duplicated rows from AddUser.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("AddUser")

# New sheet goes right after AddUser (last sheet in the workbook today).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dst = $wb.Worksheets.Add($null, $lastSheet)
$dst.Name = "FindNewUser"

# Row 1: same Section/Path/Action/Value header used on every other sheet -
# copy it (with formatting) from AddUser so the bold style matches exactly.
$src.Range("A1:D1").Copy($dst.Range("A1"))

# Rows 2-4: the old "Search" / "Enter_Text" / "Retrieve_Row_Test" rows
# (previously AddUser!A23:D25) become the body of the new sheet.
$src.Range("A23:D25").Copy($dst.Range("A2"))

# The source rows never populated B/D on the last row - Copy() leaves
# behind blank-but-styled cells there, so drop them to match the source.
$dst.Range("B4").ClearContents()
$dst.Range("D4").ClearContents()

# Match the rest of the workbook's column sizing on the new sheet.
$dst.Range("A1:D1").ColumnWidth = 14.6

# The rows have been relocated, so remove them from AddUser ...
$src.Rows("23:25").Delete()
# ... leaving the selection on the row the deleted rows used to occupy.
$src.Rows("23:23").Select() | Out-Null

# Land on the new sheet, matching its last saved selection.
$dst.Activate()
$dst.Range("D20").Select() | Out-Null
